# Updates cryptocurrency price/volume data to match latest snapshot.
# Rows 27 and 28 (Stellar/Cosmos) also swap their Coin/Link/Price/Volume data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ordered map of cell reference -> new text value.
$updates = [ordered]@{
    'D2' = '26.948.59'
    'E2' = '  +0.43%  '
    'D3' = '1.640.37'
    'E3' = '  +0.02%  '
    'E4' = '  -0.56%  '
    'D5' = '217.93'
    'E5' = '  +0.17%  '
    'D6' = '0.507'
    'E6' = '  +1.93%  '
    'E7' = '  -0.48%  '
    'E8' = '  +1.68%  '
    'E9' = '  +0.43%  '
    'D10' = '19.99'
    'E10' = '  +3.97%  '
    'E11' = '  -0.02%  '
    'D12' = '1.869.78'
    'E12' = '  -0.02%  '
    'D13' = '1.641.98'
    'E13' = '  +0.27%  '
    'E14' = '  -0.83%  '
    'E15' = '  +1.36%  '
    'D16' = '67.16'
    'E16' = '  +3.12%  '
    'D17' = '26.936.04'
    'E17' = '  +0.37%  '
    'D18' = '0.0₃0732'
    'E18' = '  +0.47%  '
    'D19' = '219.29'
    'E19' = '  +2.02%  '
    'E20' = '  -0.49%  '
    'D21' = '6.74'
    'E21' = '  +2.63%  '
    'E22' = '  +1.38%  '
    'E23' = '  +1.39%  '
    'E24' = '  +0.24%  '
    'D25' = '147.27'
    'E25' = '  +0.03%  '
    'E26' = '  -0.62%  '
    'B27' = 'Cosmos'
    'C27' = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
    'D27' = '7.31'
    'E27' = '  +1.72%  '
    'B28' = 'Stellar'
    'C28' = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
    'D28' = '0.119'
    'E28' = '  +0.98%  '
    'E29' = '  +0.17%  '
    'E30' = '  -0.71%  '
    'E31' = '  -0.41%  '
    'E32' = '  -0.59%  '
    'E33' = '  +0.60%  '
    'E34' = '  +1.01%  '
    'D35' = '1.268.93'
    'E35' = '  -0.45%  '
    'E36' = '  -0.04%  '
    'E37' = '  +2.73%  '
    'E38' = '  +1.83%  '
    'E39' = '  +2.60%  '
    'E40' = '  -0.40%  '
    'D41' = '0.807'
    'E41' = '  +0.57%  '
    'D42' = '5.35'
    'E42' = '  +0.78%  '
    'D43' = '1.780.12'
    'E43' = '  -0.05%  '
    'E44' = '  +1.25%  '
    'D45' = '62.21'
    'D46' = '92.33'
    'E46' = '  -0.11%  '
    'E47' = '  +1.45%  '
    'D48' = '0.0₆0106'
    'E48' = '  +18.35%  '
    'E49' = '  -0.75%  '
    'E50' = '  +1.61%  '
    'D51' = '0.0964'
    'E51' = '  -0.39%  '
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    # Force text storage so numeric-looking strings (e.g. '217.93')
    # and percentages keep their exact original text representation
    # instead of being auto-converted to Excel numbers.
    $cell.NumberFormat = '@'
    $cell.Value = $updates[$ref]
    $cell.Style = 'Normal'
}

Write-Output "Applied $($updates.Count) cell updates"
